$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 24 (2025-11 stats) with the refreshed figures.
$ws.Range("B24").Value = 6352
$ws.Range("C24").Value = 1000
$ws.Range("D24").Value = 5951712
$ws.Range("E24").Value = 936.9823677581863
$ws.Range("F24").Value = 8.285032390044321
$ws.Range("G24").Value = 3.626943005181338
$ws.Range("H24").Value = 26.08031399584627
